$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$lastRow = $used.Rows.Count
$lastCol = $used.Columns.Count

# Locate the "Date" header column within row 1 (robust to column drift).
$dateCol = 0
for ($c = 1; $c -le $lastCol; $c++) {
    if ($ws.Cells.Item(1, $c).Value2 -eq "Date") {
        $dateCol = $c
        break
    }
}
if ($dateCol -eq 0) {
    $dateCol = 58  # fallback: column BF
}

$rng = $ws.Range($ws.Cells.Item(2, $dateCol), $ws.Cells.Item($lastRow, $dateCol))
# Force text interpretation so the ISO-formatted replacement date string
# isn't auto-converted into a date serial value by Excel's input parser.
$rng.NumberFormat = "@"

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, $dateCol)
    if ($cell.Value2 -eq "4-29-2007-08") {
        $cell.Value = "2008-04-29"
    }
}

# Restore the default cell style so these cells keep the workbook's
# original (unstyled) appearance rather than an explicit text format.
$rng.Style = "Normal"
